# Update cryptos price (D) and volume-change (E) columns per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.541.27"
$ws.Range("E2").Value = "  -3.91%  "
$ws.Range("D3").Value = "2.508.48"
$ws.Range("E3").Value = "  -5.02%  "
$ws.Range("D5").Value = "578.61"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "166.86"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "2.505.94"
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("E10").Value = "  -7.16%  "
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "2.966.75"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").Value = "69.394.15"
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("E16").Value = "  -6.15%  "
$ws.Range("D17").Value = "24.95"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").Value = "2.533.52"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("E19").Value = "  -7.22%  "
$ws.Range("D20").Value = "11.36"
$ws.Range("E20").Value = "  -6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.61%  "
$ws.Range("D22").Value = "3.95"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").Value = "3.99"
$ws.Range("E26").Value = "  -5.91%  "
$ws.Range("D27").Value = "8.93"
$ws.Range("E27").Value = "  -5.74%  "
$ws.Range("D28").Value = "2.635.05"
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "462.94"
$ws.Range("E33").Value = "  -6.50%  "
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +3.74%  "
$ws.Range("D37").Value = "154.39"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").Value = "18.96"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "18.41"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "4.76"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  -6.72%  "
$ws.Range("D44").Value = "1.17"
$ws.Range("E44").Value = "  -14.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.37%  "
$ws.Range("D46").Value = "38.11"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "143.24"
$ws.Range("E47").Value = "  -5.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("E51").Value = "  -1.99%  "
